$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1: copy style from E1 (bold header formatting), then set text
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Data cells F2:F69: time_taken timestamps
$ws.Cells.Item(2, 6).Value = "2021-10-05 13:42:15.222855"
$ws.Cells.Item(3, 6).Value = "2021-10-05 13:42:15.222868"
$ws.Cells.Item(4, 6).Value = "2021-10-05 13:42:15.222872"
$ws.Cells.Item(5, 6).Value = "2021-10-05 13:42:15.222876"
$ws.Cells.Item(6, 6).Value = "2021-10-05 13:42:15.222879"
$ws.Cells.Item(7, 6).Value = "2021-10-05 13:42:15.222883"
$ws.Cells.Item(8, 6).Value = "2021-10-05 13:42:15.222886"
$ws.Cells.Item(9, 6).Value = "2021-10-05 13:42:15.222888"
$ws.Cells.Item(10, 6).Value = "2021-10-05 13:42:15.222892"
$ws.Cells.Item(11, 6).Value = "2021-10-05 13:42:15.222895"
$ws.Cells.Item(12, 6).Value = "2021-10-05 13:42:15.222898"
$ws.Cells.Item(13, 6).Value = "2021-10-05 13:42:15.222901"
$ws.Cells.Item(14, 6).Value = "2021-10-05 13:42:15.222904"
$ws.Cells.Item(15, 6).Value = "2021-10-05 13:42:15.222907"
$ws.Cells.Item(16, 6).Value = "2021-10-05 13:42:15.222910"
$ws.Cells.Item(17, 6).Value = "2021-10-05 13:42:15.222913"
$ws.Cells.Item(18, 6).Value = "2021-10-05 13:42:15.222916"
$ws.Cells.Item(19, 6).Value = "2021-10-05 13:42:15.222920"
$ws.Cells.Item(20, 6).Value = "2021-10-05 13:42:15.222923"
$ws.Cells.Item(21, 6).Value = "2021-10-05 13:42:15.222926"
$ws.Cells.Item(22, 6).Value = "2021-10-05 13:42:15.222929"
$ws.Cells.Item(23, 6).Value = "2021-10-05 13:42:15.222932"
$ws.Cells.Item(24, 6).Value = "2021-10-05 13:42:15.222935"
$ws.Cells.Item(25, 6).Value = "2021-10-05 13:42:15.222937"
$ws.Cells.Item(26, 6).Value = "2021-10-05 13:42:15.222941"
$ws.Cells.Item(27, 6).Value = "2021-10-05 13:42:15.222944"
$ws.Cells.Item(28, 6).Value = "2021-10-05 13:42:15.222947"
$ws.Cells.Item(29, 6).Value = "2021-10-05 13:42:15.222950"
$ws.Cells.Item(30, 6).Value = "2021-10-05 13:42:15.222953"
$ws.Cells.Item(31, 6).Value = "2021-10-05 13:42:15.222956"
$ws.Cells.Item(32, 6).Value = "2021-10-05 13:42:15.222959"
$ws.Cells.Item(33, 6).Value = "2021-10-05 13:42:15.222962"
$ws.Cells.Item(34, 6).Value = "2021-10-05 13:42:15.222966"
$ws.Cells.Item(35, 6).Value = "2021-10-05 13:42:15.222969"
$ws.Cells.Item(36, 6).Value = "2021-10-05 13:42:15.222972"
$ws.Cells.Item(37, 6).Value = "2021-10-05 13:42:15.222975"
$ws.Cells.Item(38, 6).Value = "2021-10-05 13:42:15.222978"
$ws.Cells.Item(39, 6).Value = "2021-10-05 13:42:15.222980"
$ws.Cells.Item(40, 6).Value = "2021-10-05 13:42:15.222983"
$ws.Cells.Item(41, 6).Value = "2021-10-05 13:42:15.222986"
$ws.Cells.Item(42, 6).Value = "2021-10-05 13:42:15.222990"
$ws.Cells.Item(43, 6).Value = "2021-10-05 13:42:15.222993"
$ws.Cells.Item(44, 6).Value = "2021-10-05 13:42:15.222996"
$ws.Cells.Item(45, 6).Value = "2021-10-05 13:42:15.222999"
$ws.Cells.Item(46, 6).Value = "2021-10-05 13:42:15.223002"
$ws.Cells.Item(47, 6).Value = "2021-10-05 13:42:15.223005"
$ws.Cells.Item(48, 6).Value = "2021-10-05 13:42:15.223008"
$ws.Cells.Item(49, 6).Value = "2021-10-05 13:42:15.223011"
$ws.Cells.Item(50, 6).Value = "2021-10-05 13:42:15.223014"
$ws.Cells.Item(51, 6).Value = "2021-10-05 13:42:15.223017"
$ws.Cells.Item(52, 6).Value = "2021-10-05 13:42:15.223019"
$ws.Cells.Item(53, 6).Value = "2021-10-05 13:42:15.223022"
$ws.Cells.Item(54, 6).Value = "2021-10-05 13:42:15.223026"
$ws.Cells.Item(55, 6).Value = "2021-10-05 13:42:15.223029"
$ws.Cells.Item(56, 6).Value = "2021-10-05 13:42:15.223032"
$ws.Cells.Item(57, 6).Value = "2021-10-05 13:42:15.223035"
$ws.Cells.Item(58, 6).Value = "2021-10-05 13:42:15.223038"
$ws.Cells.Item(59, 6).Value = "2021-10-05 13:42:15.223041"
$ws.Cells.Item(60, 6).Value = "2021-10-05 13:42:15.223044"
$ws.Cells.Item(61, 6).Value = "2021-10-05 13:42:15.223047"
$ws.Cells.Item(62, 6).Value = "2021-10-05 13:42:15.223050"
$ws.Cells.Item(63, 6).Value = "2021-10-05 13:42:15.223053"
$ws.Cells.Item(64, 6).Value = "2021-10-05 13:42:15.223056"
$ws.Cells.Item(65, 6).Value = "2021-10-05 13:42:15.223059"
$ws.Cells.Item(66, 6).Value = "2021-10-05 13:42:15.223064"
$ws.Cells.Item(67, 6).Value = "2021-10-05 13:42:15.223067"
$ws.Cells.Item(68, 6).Value = "2021-10-05 13:42:15.223070"
$ws.Cells.Item(69, 6).Value = "2021-10-05 13:42:15.223073"
